# Fix the Excel template bug: the "Tanque" sheet's validation messages
# in column F were hard-coded to "Validado com sucesso! ..." instead of
# reporting the actual SPED vs. relatório divergence values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tanque")

$ws.Range("F2").Value = "Divergência entre o SPED(7526,00) e o relatório(6800,75)!"
$ws.Range("F3").Value = "Divergência entre o SPED(7526,00) e o relatório(7301,36)!"
$ws.Range("F4").Value = "Divergência entre o SPED(7526,00) e o relatório(5913,61)!"
$ws.Range("F5").Value = "Divergência entre o SPED(7526,00) e o relatório(6713,67)!"
